# Docs files and other.
#
# Rename the two worksheets, update the formula/data on the BLEU sheet,
# widen its third column, and move the "active" tab / cursor position
# from the CompilationAndTestsRan sheet to the BLEU sheet.

$wb = $excel.ActiveWorkbook

$wsBleu = $wb.Worksheets.Item(1)
$wsComp = $wb.Worksheets.Item(2)

# --- Rename sheets (defined names referencing them are updated automatically) ---
$wsBleu.Name = "BLEUBasic"
$wsComp.Name = "CompilationBasic"

# --- Touch the CompilationBasic sheet's view first (so it is no longer the
#     active tab at the end) while preserving its existing selection state ---
$wsComp.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1

# --- Add the new total row under the BLEU data and set the new column width ---
$wsBleu.Activate()
$wsBleu.Range("C171").Formula = "=SUM(C2:C170)"
$wsBleu.Columns.Item(3).ColumnWidth = 32.5

# --- Leave the BLEU sheet as the active tab, with C171 selected and the
#     window scrolled down to row 112, matching the new view state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 112
$win.ScrollColumn = 1
[void]$wsBleu.Range("C171").Select()
